$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2, E2 and F2 would otherwise be auto-coerced by Excel into a number / date
# serial. Briefly force a text format so the literal string is preserved,
# then clear the format again so no stray style survives on the cell.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "400011172559639"
$ws.Range("A2").ClearFormats()

$ws.Range("D2").Value = "12Μ0ΤΔΑ"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "10768"
$ws.Range("E2").ClearFormats()

$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "03/10/2025"
$ws.Range("F2").ClearFormats()

$ws.Range("I2").Value = "13,83"
$ws.Range("J2").Value = "3,32"
$ws.Range("K2").Value = "17,15"
